$wb = $excel.ActiveWorkbook

# Sheet "2025"
$ws = $wb.Worksheets.Item("2025")
$ws.Range("N2").Value = 7155.07680048089
$ws.Range("O2").Value = 6980.325837388836

# Sheet "2030"
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 5707.815717280662
$ws.Range("I2").Value = 44492.05901988943
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 21991.42050229464
$ws.Range("O2").Value = 12079.4099131153

# Sheet "2035"
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 2927.360317916481
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15117.91314604427
$ws.Range("O2").Value = 14761.05508568936

# Sheet "2040"
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 2927.360317916481
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15225.03710440102
$ws.Range("O2").Value = 14761.05508568936

# Sheet "2045"
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 6352.985609279765
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15769.76495484199
$ws.Range("O2").Value = 17096.52172347162

# Sheet "2050"
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 6352.985609279765
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15769.76495484199
$ws.Range("O2").Value = 17096.52172347162
